# Updates the cryptos list (Price and Volume(1h) columns) with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.539.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "'1.628.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D5").Value = "'212.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'19.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "'1.855.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'1.631.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'4.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "'63.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'26.582.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "'214.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.33%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'4.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").Value = "'148.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").Value = "'15.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "'0.0505"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "'1.219.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "'0.0174"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.37%  "
$ws.Range("D39").Value = "'0.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'2.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").Value = "'0.794"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").Value = "'1.766.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'92.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'54.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'7.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "'0.409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
